$wb = $excel.ActiveWorkbook

# Update Yufeng Lai's attendance numbers on the "2024 - Fall" sheet.
$wsFall = $wb.Worksheets.Item("2024 - Fall")
$wsFall.Range("G3").Value = 18
$wsFall.Range("H3").Value = 4

# Move the selection on the "2024 - Fall" sheet to G4 (matches author's recorded view state).
$wsFall.Range("G4").Select()

# Make "Attendance Descriptives" the active sheet/tab, keeping its prior selection (C15).
$wsDesc = $wb.Worksheets.Item("Attendance Descriptives")
$wsDesc.Activate()
$wsDesc.Range("C15").Select()
